$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1283.3334
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 1283.3334
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = -1113.3334
$ws.Range("N12").Value = -490

# Row 15 (Leve Item ID 44146)
$ws.Range("H15").Value = 405
$ws.Range("I15").Value = 405
$ws.Range("K15").Value = 1215
$ws.Range("M15").Value = -1046

# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 1805.4546
$ws.Range("J17").Value = 1805.4546
$ws.Range("L17").Value = 5416.3638
$ws.Range("N17").Value = -5752.3638

# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 10905
$ws.Range("I18").Value = 10255.357
$ws.Range("K18").Value = 10255.357
$ws.Range("M18").Value = -9971.357

# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 10999.75
$ws.Range("J32").Value = 10999.75
$ws.Range("L32").Value = 10999.75
$ws.Range("N32").Value = -11651.75

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 6249.25
$ws.Range("I43").Value = 6249.25
$ws.Range("K43").Value = 6249.25
$ws.Range("M43").Value = -6180.25

# Row 109 (Leve Item ID 25639)
$ws.Range("H109").Value = 85497
$ws.Range("J109").Value = 85497
$ws.Range("L109").Value = 85497
$ws.Range("N109").Value = -88271

# Row 114 (Leve Item ID 25959)
$ws.Range("H114").Value = 94997.5
$ws.Range("J114").Value = 99995
$ws.Range("L114").Value = 99995
$ws.Range("N114").Value = -108673

# Row 117 (Leve Item ID 26118)
$ws.Range("H117").Value = 99995
$ws.Range("J117").Value = 99995
$ws.Range("L117").Value = 99995
$ws.Range("N117").Value = -109173

# Row 120 (Leve Item ID 26279)
$ws.Range("H120").Value = 99995
$ws.Range("J120").Value = 99995
$ws.Range("L120").Value = 99995
$ws.Range("N120").Value = -109671

# Row 124 (Leve Item ID 34241)
$ws.Range("H124").Value = 99995
$ws.Range("J124").Value = 99995
$ws.Range("L124").Value = 99995
$ws.Range("N124").Value = -109815

# Row 128 (Leve Item ID 34540)
$ws.Range("H128").Value = 99995
$ws.Range("J128").Value = 99995
$ws.Range("L128").Value = 99995
$ws.Range("N128").Value = -109955

# Row 130 (Leve Item ID 34691)
$ws.Range("H130").Value = 76990.39999999999
$ws.Range("J130").Value = 76990.39999999999
$ws.Range("L130").Value = 76990.39999999999
$ws.Range("N130").Value = -87030.39999999999

# Row 133 (Leve Item ID 41856)
$ws.Range("H133").Value = 99995
$ws.Range("J133").Value = 99995
$ws.Range("L133").Value = 99995
$ws.Range("N133").Value = -110115

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2235.2942
$ws.Range("I137").Value = 2500
$ws.Range("J137").Value = 2200
$ws.Range("K137").Value = 7500
$ws.Range("L137").Value = 6600
$ws.Range("M137").Value = -4950
$ws.Range("N137").Value = -11700

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 7749.25
$ws.Range("J138").Value = 7749.25
$ws.Range("L138").Value = 23247.75
$ws.Range("N138").Value = -33527.75

# Row 140 (Leve Item ID 42459)
$ws.Range("H140").Value = 79997
$ws.Range("J140").Value = 79997
$ws.Range("L140").Value = 79997
$ws.Range("N140").Value = -90357

$ws = $wb.Worksheets.Item("ARM")
# Row 17 (Leve Item ID 2495)
$ws.Range("H17").Value = 903.2
$ws.Range("J17").Value = 903.2
$ws.Range("L17").Value = 903.2
$ws.Range("N17").Value = -1249.2

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 415
$ws.Range("I61").Value = 415
$ws.Range("K61").Value = 415
$ws.Range("M61").Value = -203

# Row 95 (Leve Item ID 18204)
$ws.Range("H95").Value = 2000000
$ws.Range("J95").Value = 2000000
$ws.Range("L95").Value = 2000000
$ws.Range("N95").Value = -2005492

# Row 113 (Leve Item ID 26002)
$ws.Range("H113").Value = 40198.5
$ws.Range("J113").Value = 40198.5
$ws.Range("L113").Value = 40198.5
$ws.Range("N113").Value = -48876.5

# Row 119 (Leve Item ID 26287)
$ws.Range("H119").Value = 698
$ws.Range("J119").Value = 698
$ws.Range("L119").Value = 698
$ws.Range("N119").Value = -10374

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122:N122").ClearContents()

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 415
$ws.Range("I136").Value = 415
$ws.Range("K136").Value = 1245
$ws.Range("M136").Value = 1305

# Row 140 (Leve Item ID 42496)
$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 3231.2
$ws.Range("I107").Value = 1385.3334
$ws.Range("K107").Value = 1385.3334
$ws.Range("M107").Value = 534.6666

# Row 110 (Leve Item ID 25790)
$ws.Range("H110").Value = 40702
$ws.Range("J110").Value = 40702
$ws.Range("L110").Value = 40702
$ws.Range("N110").Value = -48882

# Row 111 (Leve Item ID 25789)
$ws.Range("H111").Value = 98997
$ws.Range("J111").Value = 98997
$ws.Range("L111").Value = 98997
$ws.Range("N111").Value = -107177

# Row 120 (Leve Item ID 26275)
$ws.Range("H120").Value = 99995
$ws.Range("J120").Value = 99995
$ws.Range("L120").Value = 99995
$ws.Range("N120").Value = -109671

# Row 130 (Leve Item ID 34682)
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 132 (Leve Item ID 41855)
$ws.Range("H132").Value = 99000
$ws.Range("J132").Value = 99000
$ws.Range("L132").Value = 99000
$ws.Range("N132").Value = -109120

# Row 140 (Leve Item ID 42471)
$ws.Range("H140").Value = 90680
$ws.Range("J140").Value = 90680
$ws.Range("L140").Value = 90680
$ws.Range("N140").Value = -101040

$ws = $wb.Worksheets.Item("CRP")
# Row 20 (Leve Item ID 34533)
$ws.Range("H20").Value = 94500
$ws.Range("J20").Value = 94500
$ws.Range("L20").Value = 94500
$ws.Range("N20").Value = -94972

# Row 30 (Leve Item ID 34533)
$ws.Range("H30").Value = 94500
$ws.Range("J30").Value = 94500
$ws.Range("L30").Value = 94500
$ws.Range("N30").Value = -94682

# Row 36 (Leve Item ID 1845)
$ws.Range("H36").Value = 4048
$ws.Range("I36").Value = 4048
$ws.Range("K36").Value = 4048
$ws.Range("M36").Value = -3660

# Row 40 (Leve Item ID 1845)
$ws.Range("H40").Value = 4048
$ws.Range("I40").Value = 4048
$ws.Range("K40").Value = 4048
$ws.Range("M40").Value = -3888

# Row 112 (Leve Item ID 25796)
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# Row 118 (Leve Item ID 26112)
$ws.Range("H118").Value = 79997
$ws.Range("J118").Value = 79997
$ws.Range("L118").Value = 79997
$ws.Range("N118").Value = -83311

# Row 119 (Leve Item ID 26276)
$ws.Range("H119").Value = 99995
$ws.Range("J119").Value = 99995
$ws.Range("L119").Value = 99995
$ws.Range("N119").Value = -109671

# Row 128 (Leve Item ID 34533)
$ws.Range("H128").Value = 94500
$ws.Range("J128").Value = 94500
$ws.Range("L128").Value = 94500
$ws.Range("N128").Value = -104460

# Row 130 (Leve Item ID 34689)
$ws.Range("H130").Value = 99995
$ws.Range("J130").Value = 99995
$ws.Range("L130").Value = 99995
$ws.Range("N130").Value = -110035

# Row 135 (Leve Item ID 42008)
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140

$ws = $wb.Worksheets.Item("CUL")
# Row 6 (Leve Item ID 4639)
$ws.Range("H6").Value = 49.5
$ws.Range("I6").Value = 54.375
$ws.Range("K6").Value = 163.125
$ws.Range("M6").Value = -50.125

# Row 12 (Leve Item ID 4854)
$ws.Range("H12").Value = 21
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 21
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 63
$ws.Range("N12").Value = -409
$ws.Range("M12").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 110 (Leve Item ID 25802)
$ws.Range("H110").Value = 99995
$ws.Range("J110").Value = 99995
$ws.Range("L110").Value = 99995
$ws.Range("N110").Value = -108175

# Row 116 (Leve Item ID 26120)
$ws.Range("H116").Value = 99995
$ws.Range("J116").Value = 99995
$ws.Range("L116").Value = 99995
$ws.Range("N116").Value = -109173

# Row 119 (Leve Item ID 26282)
$ws.Range("H119").Value = 99995
$ws.Range("J119").Value = 99995
$ws.Range("L119").Value = 99995
$ws.Range("N119").Value = -109671

# Row 124 (Leve Item ID 34247)
$ws.Range("H124").Value = 99995
$ws.Range("J124").Value = 99995
$ws.Range("L124").Value = 99995
$ws.Range("N124").Value = -109815

# Row 128 (Leve Item ID 34544)
$ws.Range("H128").Value = 99995
$ws.Range("J128").Value = 99995
$ws.Range("L128").Value = 99995
$ws.Range("N128").Value = -109955

# Row 130 (Leve Item ID 34692)
$ws.Range("H130").Value = 99995
$ws.Range("J130").Value = 99995
$ws.Range("L130").Value = 99995
$ws.Range("N130").Value = -110035

# Row 133 (Leve Item ID 41854)
$ws.Range("H133").Value = 94997.5
$ws.Range("J133").Value = 94997.5
$ws.Range("L133").Value = 94997.5
$ws.Range("N133").Value = -105117.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 4072
$ws.Range("J46").Value = 4072
$ws.Range("L46").Value = 4072
$ws.Range("N46").Value = -4448

# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 2650.25
$ws.Range("I61").Value = 2650.25
$ws.Range("K61").Value = 2650.25
$ws.Range("M61").Value = -2448.25

# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 2650.25
$ws.Range("I113").Value = 2650.25
$ws.Range("K113").Value = 2650.25
$ws.Range("M113").Value = -480.25

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 5574.7
$ws.Range("I96").Value = 4041.1667
$ws.Range("J96").Value = 7875
$ws.Range("K96").Value = 4041.1667
$ws.Range("L96").Value = 7875
$ws.Range("M96").Value = -2668.1667
$ws.Range("N96").Value = -10621

# Row 97 (Leve Item ID 18220)
$ws.Range("H97").Value = 42665.668
$ws.Range("J97").Value = 42665.668
$ws.Range("L97").Value = 42665.668
$ws.Range("N97").Value = -44647.668

# Row 140 (Leve Item ID 42506)
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
